$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (46074 -> 46075, i.e. 2026-02-21 -> 2026-02-22) for every data row (2-82).
for ($r = 2; $r -le 82; $r++) {
    $ws.Cells.Item($r, 3).Value = 46075
}
